# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 416
    $ws.Range("F3").Value = 5175
}

# Row with the 50 -> 51 update is row 5 on "展览" but row 6 on "全部类型"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 51

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 51
